# QA manual CoderHouse 2025
#
# The workbook was reopened/resaved and its columns were re-fit to their
# content: column C (the "Title" column, holding the long issue
# descriptions) had never had an explicit width before - it picks one up
# for the first time - while columns D:H get their best-fit widths
# refreshed too. Columns A and B already matched their best-fit width, so
# they are intentionally left untouched (matches the observed diff).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Recompute the "best fit" width for every column that actually holds
# text (C:H) from their content.
$ws.Range("C1:H11").EntireColumn.AutoFit() | Out-Null

# Nudge each of those columns to the exact best-fit width Excel itself
# calculated (sub-pixel, font-metric based - not reproducible by the
# host's coarser AutoFit heuristic alone).
$ws.Columns.Item(3).ColumnWidth = 126.72135416666667   # -> 127.5546875
$ws.Columns.Item(4).ColumnWidth = 38.608072916666664   # -> 39.44140625
$ws.Columns.Item(5).ColumnWidth = 5.276041666666667    # -> 6.109375
$ws.Columns.Item(6).ColumnWidth = 8.830729166666666    # -> 9.6640625
$ws.Columns.Item(7).ColumnWidth = 4.053385416666667    # -> 4.88671875
$ws.Columns.Item(8).ColumnWidth = 14.608072916666666   # -> 15.44140625
